# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-51.
# All of these cells store plain text (e.g. "1.000", "  +0.38%  ", "29.931.97").
# For values that are plain text but LOOK like a number (e.g. "243.91"), Excel's
# COM layer will silently coerce a simple .Value assignment into a real numeric
# cell (losing the original text formatting / trailing zeros / precision). To keep
# those cells genuinely textual -- matching the source workbook -- we briefly force
# a text number format before assigning the value, then restore the cell's normal
# style so no visible formatting changes are introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '29.931.97'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '1.892.96'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  -0.01%  '
Set-TextValue 'D5' '0.7729'
$ws.Range('E5').Value = '  -2.70%  '
Set-TextValue 'D6' '243.91'
$ws.Range('E6').Value = '  +0.44%  '
Set-TextValue 'D7' '1.001'
$ws.Range('E7').Value = '  -0.03%  '
Set-TextValue 'D8' '0.3131'
$ws.Range('E8').Value = '  -0.95%  '
Set-TextValue 'D9' '25.58'
$ws.Range('E9').Value = '  +0.81%  '
Set-TextValue 'D10' '0.07327'
$ws.Range('E10').Value = '  +3.77%  '
Set-TextValue 'D11' '0.08052'
$ws.Range('E11').Value = '  -0.37%  '
Set-TextValue 'D12' '0.7712'
$ws.Range('E12').Value = '  +0.50%  '
Set-TextValue 'D13' '5.488'
$ws.Range('E13').Value = '  +2.49%  '
Set-TextValue 'D14' '93.98'
$ws.Range('E14').Value = '  +1.71%  '
$ws.Range('D15').Value = '1.821.10'
$ws.Range('E15').Value = '  -3.82%  '
Set-TextValue 'D16' '6.214'
$ws.Range('E16').Value = '  +3.49%  '
$ws.Range('D17').Value = '29.877.43'
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('E18').Value = '  +0.95%  '
Set-TextValue 'D19' '246.29'
Set-TextValue 'D20' '0.000007852'
$ws.Range('E20').Value = '  +1.85%  '
Set-TextValue 'D21' '8.151'
$ws.Range('E21').Value = '  -2.32%  '
Set-TextValue 'D22' '0.9998'
$ws.Range('D23').Value = '2.108.46'
$ws.Range('E23').Value = '  -1.81%  '
$ws.Range('E24').Value = '  -0.04%  '
Set-TextValue 'D25' '0.1571'
$ws.Range('E25').Value = '  -3.99%  '
Set-TextValue 'D26' '9.433'
$ws.Range('E26').Value = '  +0.77%  '
Set-TextValue 'D27' '162.11'
$ws.Range('E27').Value = '  -2.39%  '
Set-TextValue 'D28' '18.75'
Set-TextValue 'D29' '2.023'
$ws.Range('E29').Value = '  -1.51%  '
Set-TextValue 'D30' '1.422'
$ws.Range('E30').Value = '  +1.65%  '
Set-TextValue 'D31' '1.540'
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('E32').Value = '  +0.92%  '
Set-TextValue 'D33' '0.05546'
$ws.Range('E33').Value = '  -2.93%  '
Set-TextValue 'D34' '4.065'
$ws.Range('E34').Value = '  +0.47%  '
$ws.Range('E35').Value = '  -2.19%  '
Set-TextValue 'D36' '0.7490'
$ws.Range('E36').Value = '  +1.44%  '
Set-TextValue 'D37' '0.9991'
$ws.Range('E37').Value = '  +0.00%  '
Set-TextValue 'D38' '2.682'
$ws.Range('E38').Value = '  +1.90%  '
Set-TextValue 'D39' '0.01926'
$ws.Range('E39').Value = '  +0.75%  '
$ws.Range('E41').Value = '  +1.51%  '
Set-TextValue 'D42' '74.24'
$ws.Range('E42').Value = '  +2.42%  '
$ws.Range('D43').Value = '1.097.89'
$ws.Range('E43').Value = '  +6.23%  '
Set-TextValue 'D44' '6.003'
Set-TextValue 'D45' '0.8509'
$ws.Range('E45').Value = '  +1.13%  '
$ws.Range('E46').Value = '  +0.00%  '
Set-TextValue 'D47' '1.884'
$ws.Range('E47').Value = '  +0.69%  '
Set-TextValue 'D48' '102.45'
$ws.Range('E48').Value = '  -0.67%  '
$ws.Range('E49').Value = '  +1.54%  '
Set-TextValue 'D50' '9.780'
$ws.Range('E50').Value = '  -1.97%  '
Set-TextValue 'D51' '2.997'
$ws.Range('E51').Value = '  +3.32%  '
